$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Teste Covid"
$ws.Range("B1").Value = 500
$ws.Range("C1").Value = "V. Epidemiológica"
$ws.Range("D1").Value = "Teste Covid"

$ws.Range("A2").Value = "Alfinete"
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = "CERPAT"
$ws.Range("D2").Value = "Teste Covid"

$ws.Range("A3").Value = "Papel"
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = "V. Sanitária"
$ws.Range("D3").Value = "Material"

$ws.Range("A1:F8").Select() | Out-Null
